# Swissbix template.docx edit:
# - Insert a new "Definizione Economica" heading paragraph followed by a
#   pricing table (header row + Jinja "for item in items" loop rows) right
#   after the page-break paragraph that currently ends the first page,
#   and before the existing blank paragraph (w14:paraId 6D116AFD).

$d = $word.ActiveDocument

# Locate the insertion point: end of the first paragraph in the document
# body (the one holding the trailing page breaks / single space run).
$insertPos = $d.Paragraphs.Item(1).Range.End
$insertionRange = $d.Range($insertPos, $insertPos)

$fragmentXml = '<w:p><w:r><w:rPr><w:rFonts w:ascii="Calibri" w:hAnsi="Calibri"/><w:color w:val="C00000"/><w:sz w:val="30"/></w:rPr><w:t>Definizione Economica</w:t></w:r></w:p><w:tbl><w:tblPr><w:tblStyle w:val="bixstyle"/><w:tblW w:w="0" w:type="auto"/><w:tblLook w:val="04A0" w:firstRow="1" w:lastRow="0" w:firstColumn="1" w:lastColumn="0" w:noHBand="0" w:noVBand="1"/></w:tblPr><w:tblGrid><w:gridCol w:w="2333"/><w:gridCol w:w="2333"/><w:gridCol w:w="2333"/><w:gridCol w:w="2333"/></w:tblGrid><w:tr><w:trPr><w:cnfStyle w:val="100000000000" w:firstRow="1" w:lastRow="0" w:firstColumn="0" w:lastColumn="0" w:oddVBand="0" w:evenVBand="0" w:oddHBand="0" w:evenHBand="0" w:firstRowFirstColumn="0" w:firstRowLastColumn="0" w:lastRowFirstColumn="0" w:lastRowLastColumn="0"/></w:trPr><w:tc><w:tcPr><w:tcW w:w="2333" w:type="dxa"/></w:tcPr><w:p><w:r><w:t>Descrizione</w:t></w:r></w:p></w:tc><w:tc><w:tcPr><w:tcW w:w="2333" w:type="dxa"/></w:tcPr><w:p><w:r><w:t>Qt.</w:t></w:r></w:p></w:tc><w:tc><w:tcPr><w:tcW w:w="2333" w:type="dxa"/></w:tcPr><w:p><w:r><w:t>Prezzo unitario</w:t></w:r></w:p></w:tc><w:tc><w:tcPr><w:tcW w:w="2333" w:type="dxa"/></w:tcPr><w:p><w:r><w:t>Prezzo totale</w:t></w:r></w:p></w:tc></w:tr><w:tr><w:tc><w:tcPr><w:tcW w:w="9332" w:type="dxa"/><w:gridSpan w:val="4"/></w:tcPr><w:p><w:r><w:t>{%</w:t></w:r><w:proofErr w:type="spellStart"/><w:r><w:t>tr</w:t></w:r><w:proofErr w:type="spellEnd"/><w:r><w:t xml:space="preserve"> for item in items %}</w:t></w:r></w:p></w:tc></w:tr><w:tr><w:tc><w:tcPr><w:tcW w:w="2333" w:type="dxa"/></w:tcPr><w:p><w:proofErr w:type="gramStart"/><w:r><w:t xml:space="preserve">{{ </w:t></w:r><w:proofErr w:type="spellStart"/><w:r><w:t>item</w:t></w:r><w:proofErr w:type="gramEnd"/><w:r><w:t>.descrizione</w:t></w:r><w:proofErr w:type="spellEnd"/><w:r><w:t xml:space="preserve"> }}</w:t></w:r></w:p></w:tc><w:tc><w:tcPr><w:tcW w:w="2333" w:type="dxa"/></w:tcPr><w:p><w:proofErr w:type="gramStart"/><w:r><w:t xml:space="preserve">{{ </w:t></w:r><w:proofErr w:type="spellStart"/><w:r><w:t>item.qt</w:t></w:r><w:proofErr w:type="spellEnd"/><w:proofErr w:type="gramEnd"/><w:r><w:t xml:space="preserve"> }}</w:t></w:r></w:p></w:tc><w:tc><w:tcPr><w:tcW w:w="2333" w:type="dxa"/></w:tcPr><w:p><w:proofErr w:type="gramStart"/><w:r><w:t xml:space="preserve">{{ </w:t></w:r><w:proofErr w:type="spellStart"/><w:r><w:t>item</w:t></w:r><w:proofErr w:type="gramEnd"/><w:r><w:t>.prezzo_unitario</w:t></w:r><w:proofErr w:type="spellEnd"/><w:r><w:t xml:space="preserve"> }}</w:t></w:r></w:p></w:tc><w:tc><w:tcPr><w:tcW w:w="2333" w:type="dxa"/></w:tcPr><w:p><w:proofErr w:type="gramStart"/><w:r><w:t xml:space="preserve">{{ </w:t></w:r><w:proofErr w:type="spellStart"/><w:r><w:t>item</w:t></w:r><w:proofErr w:type="gramEnd"/><w:r><w:t>.prezzo_totale</w:t></w:r><w:proofErr w:type="spellEnd"/><w:r><w:t xml:space="preserve"> }}</w:t></w:r></w:p></w:tc></w:tr><w:tr><w:tc><w:tcPr><w:tcW w:w="9332" w:type="dxa"/><w:gridSpan w:val="4"/></w:tcPr><w:p><w:r><w:t>{%</w:t></w:r><w:proofErr w:type="spellStart"/><w:r><w:t>tr</w:t></w:r><w:proofErr w:type="spellEnd"/><w:r><w:t xml:space="preserve"> </w:t></w:r><w:proofErr w:type="spellStart"/><w:r><w:t>endfor</w:t></w:r><w:proofErr w:type="spellEnd"/><w:r><w:t xml:space="preserve"> %}</w:t></w:r></w:p></w:tc></w:tr></w:tbl>'

$insertionRange.InsertXML($fragmentXml)

Write-Host "Inserted 'Definizione Economica' heading + pricing table."
